$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Sending cluster (A) and Target cluster (D) labels, and numeric columns E:T for rows 2-13
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vegfc"
$ws.Range("C2").Value = "Vipr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.401172666666667
$ws.Range("H2").Value = 10.203518
$ws.Range("I2").Value = 0.5101677883321656
$ws.Range("J2").Value = 0.5101677883321655
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1752366666666667
$ws.Range("N2").Value = 0.52571
$ws.Range("O2").Value = 0.01723504119824304
$ws.Range("P2").Value = 0.01723504119824304
$ws.Range("Q2").Value = 0.5960101608644446
$ws.Range("R2").Value = 5.364091447780001
$ws.Range("S2").Value = 0.008792762849921407
$ws.Range("T2").Value = 0.008792762849921405

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vegfc"
$ws.Range("C3").Value = "Vipr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.401172666666667
$ws.Range("H3").Value = 10.203518
$ws.Range("I3").Value = 0.5101677883321656
$ws.Range("J3").Value = 0.5101677883321655
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.550379
$ws.Range("N3").Value = 16.651137
$ws.Range("O3").Value = 0.5458960875627037
$ws.Range("P3").Value = 0.5458960875627037
$ws.Range("Q3").Value = 18.87779734444067
$ws.Range("R3").Value = 169.900176099966
$ws.Range("S3").Value = 0.2784985996510467
$ws.Range("T3").Value = 0.2784985996510467

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vegfc"
$ws.Range("C4").Value = "Vipr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.401172666666667
$ws.Range("H4").Value = 10.203518
$ws.Range("I4").Value = 0.5101677883321656
$ws.Range("J4").Value = 0.5101677883321655
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01305433333333333
$ws.Range("N4").Value = 0.039163
$ws.Range("O4").Value = 0.001283932050839421
$ws.Range("P4").Value = 0.001283932050839421
$ws.Range("Q4").Value = 0.0444000417148889
$ws.Range("R4").Value = 0.3996003754340001
$ws.Range("S4").Value = 0.0006550207747455291
$ws.Range("T4").Value = 0.000655020774745529

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vegfc"
$ws.Range("C5").Value = "Vipr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.401172666666667
$ws.Range("H5").Value = 10.203518
$ws.Range("I5").Value = 0.5101677883321656
$ws.Range("J5").Value = 0.5101677883321655
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.428794333333333
$ws.Range("N5").Value = 13.286383
$ws.Range("O5").Value = 0.4355849391882138
$ws.Range("P5").Value = 0.4355849391882139
$ws.Range("Q5").Value = 15.06309423282156
$ws.Range("R5").Value = 135.567848095394
$ws.Range("S5").Value = 0.2222214050564519
$ws.Range("T5").Value = 0.2222214050564519

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vegfc"
$ws.Range("C6").Value = "Vipr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.579868
$ws.Range("H6").Value = 7.739604
$ws.Range("I6").Value = 0.3869740471126509
$ws.Range("J6").Value = 0.3869740471126508
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1752366666666667
$ws.Range("N6").Value = 0.52571
$ws.Range("O6").Value = 0.01723504119824304
$ws.Range("P6").Value = 0.01723504119824304
$ws.Range("Q6").Value = 0.45208746876
$ws.Range("R6").Value = 4.06878721884
$ws.Range("S6").Value = 0.006669513644637379
$ws.Range("T6").Value = 0.006669513644637378

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vegfc"
$ws.Range("C7").Value = "Vipr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.579868
$ws.Range("H7").Value = 7.739604
$ws.Range("I7").Value = 0.3869740471126509
$ws.Range("J7").Value = 0.3869740471126508
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.550379
$ws.Range("N7").Value = 16.651137
$ws.Range("O7").Value = 0.5458960875627037
$ws.Range("P7").Value = 0.5458960875627037
$ws.Range("Q7").Value = 14.319245169972
$ws.Range("R7").Value = 128.873206529748
$ws.Range("S7").Value = 0.2112476183071015
$ws.Range("T7").Value = 0.2112476183071015

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vegfc"
$ws.Range("C8").Value = "Vipr2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.579868
$ws.Range("H8").Value = 7.739604
$ws.Range("I8").Value = 0.3869740471126509
$ws.Range("J8").Value = 0.3869740471126508
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01305433333333333
$ws.Range("N8").Value = 0.039163
$ws.Range("O8").Value = 0.001283932050839421
$ws.Range("P8").Value = 0.001283932050839421
$ws.Range("Q8").Value = 0.033678456828
$ws.Range("R8").Value = 0.303106111452
$ws.Range("S8").Value = 0.0004968483819309767
$ws.Range("T8").Value = 0.0004968483819309766

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vegfc"
$ws.Range("C9").Value = "Vipr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.579868
$ws.Range("H9").Value = 7.739604
$ws.Range("I9").Value = 0.3869740471126509
$ws.Range("J9").Value = 0.3869740471126508
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.428794333333333
$ws.Range("N9").Value = 13.286383
$ws.Range("O9").Value = 0.4355849391882138
$ws.Range("P9").Value = 0.4355849391882139
$ws.Range("Q9").Value = 11.425704779148
$ws.Range("R9").Value = 102.831343012332
$ws.Range("S9").Value = 0.168560066778981
$ws.Range("T9").Value = 0.168560066778981

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vegfc"
$ws.Range("C10").Value = "Vipr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6857320000000001
$ws.Range("H10").Value = 2.057196
$ws.Range("I10").Value = 0.1028581645551836
$ws.Range("J10").Value = 0.1028581645551836
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1752366666666667
$ws.Range("N10").Value = 0.52571
$ws.Range("O10").Value = 0.01723504119824304
$ws.Range("P10").Value = 0.01723504119824304
$ws.Range("Q10").Value = 0.1201653899066667
$ws.Range("R10").Value = 1.08148850916
$ws.Range("S10").Value = 0.001772764703684251
$ws.Range("T10").Value = 0.00177276470368425

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Vegfc"
$ws.Range("C11").Value = "Vipr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.6857320000000001
$ws.Range("H11").Value = 2.057196
$ws.Range("I11").Value = 0.1028581645551836
$ws.Range("J11").Value = 0.1028581645551836
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.550379
$ws.Range("N11").Value = 16.651137
$ws.Range("O11").Value = 0.5458960875627037
$ws.Range("P11").Value = 0.5458960875627037
$ws.Range("Q11").Value = 3.806072492428
$ws.Range("R11").Value = 34.254652431852
$ws.Range("S11").Value = 0.05614986960455549
$ws.Range("T11").Value = 0.05614986960455547

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Vegfc"
$ws.Range("C12").Value = "Vipr2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.6857320000000001
$ws.Range("H12").Value = 2.057196
$ws.Range("I12").Value = 0.1028581645551836
$ws.Range("J12").Value = 0.1028581645551836
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.01305433333333333
$ws.Range("N12").Value = 0.039163
$ws.Range("O12").Value = 0.001283932050839421
$ws.Range("P12").Value = 0.001283932050839421
$ws.Range("Q12").Value = 0.008951774105333335
$ws.Range("R12").Value = 0.08056596694800001
$ws.Range("S12").Value = 0.0001320628941629155
$ws.Range("T12").Value = 0.0001320628941629155

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Vegfc"
$ws.Range("C13").Value = "Vipr2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.6857320000000001
$ws.Range("H13").Value = 2.057196
$ws.Range("I13").Value = 0.1028581645551836
$ws.Range("J13").Value = 0.1028581645551836
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.428794333333333
$ws.Range("N13").Value = 13.286383
$ws.Range("O13").Value = 0.4355849391882138
$ws.Range("P13").Value = 0.4355849391882139
$ws.Range("Q13").Value = 3.036965995785334
$ws.Range("R13").Value = 27.33269396206801
$ws.Range("S13").Value = 0.04480346735278094
$ws.Range("T13").Value = 0.04480346735278093
